$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Through 2021-12-14" to "Through 2021-12-15"
$ws.Name = "Through 2021-12-15"

# Update the December label to reflect the new "through" date
$ws.Range("A13").Value = "December (through 12-15)"

# Update December row (row 13) values
$ws.Range("B13").Value = 16
$ws.Range("C13").Value = 46
$ws.Range("D13").Value = 52
$ws.Range("E13").Value = 34
$ws.Range("F13").Value = 25
$ws.Range("G13").Value = 77
$ws.Range("H13").Value = 112

# Update Total row (row 14) values
$ws.Range("B14").Value = 307
$ws.Range("C14").Value = 609
$ws.Range("D14").Value = 873
$ws.Range("E14").Value = 716
$ws.Range("F14").Value = 559
$ws.Range("G14").Value = 1341
$ws.Range("H14").Value = 1755
